# Fix total marks error on the marksheet (row 11 "Marking" and row 12 "Total").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row: Right count 5 -> 4, Wrong penalty -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Total row: Right total 85 -> 68, Wrong penalty -1 -> -2, and the score summary text
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66 / 112"
